$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.035.23"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "3.124.26"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.50"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.87"
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.122.26"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.37"
$ws.Range("E11").Value = "  +2.95%  "
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.13"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").Value = "3.634.07"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("E16").Value = "  +3.11%  "
$ws.Range("D17").Value = "64.039.78"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "3.132.58"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.84"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "480.31"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.53"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.707"
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.82"
$ws.Range("E24").Value = "  +6.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.36"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.31"
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.13"
$ws.Range("E29").Value = "  +4.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.06"
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("E31").Value = "  -6.00%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.73"
$ws.Range("E33").Value = "  +2.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.65"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("E35").Value = "  -1.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.00"
$ws.Range("E36").Value = "  +1.78%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0751"
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.61"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.98"
$ws.Range("E39").Value = "  +1.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "436.60"
$ws.Range("E40").Value = "  -3.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0391"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D44").Value = "2.861.04"
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("E45").Value = "  -1.52%  "
$ws.Range("E46").Value = "  -2.18%  "
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.81"
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("E50").Value = "  +0.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.15"
$ws.Range("E51").Value = "  +1.93%  "
